$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 454, shifting rows 454:523 down to 455:524
$ws.Rows.Item(454).Insert()

# Populate the newly inserted row 454 with the new record
$ws.Cells.Item(454, 1).Value = 6
$ws.Cells.Item(454, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(454, 3).Value = "Metropolitana"
$ws.Cells.Item(454, 4).Value = 44505
$ws.Cells.Item(454, 4).Style = $ws.Cells.Item(455, 4).Style
$ws.Cells.Item(454, 4).NumberFormat = $ws.Cells.Item(455, 4).NumberFormat
$ws.Cells.Item(454, 5).Value = 13
$ws.Cells.Item(454, 6).Value = 100112031
$ws.Cells.Item(454, 7).Value = "Poroto verde"
$ws.Cells.Item(454, 8).Value = "Magnum"
$ws.Cells.Item(454, 9).Value = "Primera"
$ws.Cells.Item(454, 10).Value = 400
$ws.Cells.Item(454, 11).Value = 38000
$ws.Cells.Item(454, 12).Value = 40000
$ws.Cells.Item(454, 13).Value = 39150
$ws.Cells.Item(454, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(454, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(454, 16).Value = 1566
$ws.Cells.Item(454, 17).Value = 25
$ws.Cells.Item(454, 18).Value = "Hortaliza"
